# Git workshop slides - minor git code and slides edits.
#
# 1) On the "Group Exercise!" slide (slide 8), insert a new numbered
#    task "Clone repo!" before the existing "Tasks" line, and turn both
#    "Clone repo!" and "Tasks" into an auto-numbered (arabic period) list.
# 2) Rename the "Eric: square" task to "Eric: exponent".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

# Insert the new first line ("Clone repo!") ahead of "Tasks", as its own
# paragraph (CR = new a:p). This must run on the shape's top-level
# TextRange so the insertion produces a real paragraph break instead of
# a literal line-break character.
$null = $tr.InsertBefore("Clone repo!" + [char]13)

# Turn "Clone repo!" into an arabic-period auto-numbered bullet.
$cloneLine = $tr.Find("Clone repo!")
$cloneBullet = $cloneLine.ParagraphFormat.Bullet
$cloneBullet.Font.Name = "+mj-lt"
$cloneBullet.Type = 2

# Turn the existing "Tasks" line into the same kind of auto-numbered bullet.
$tasksLine = $tr.Find("Tasks")
$tasksBullet = $tasksLine.ParagraphFormat.Bullet
$tasksBullet.Font.Name = "+mj-lt"
$tasksBullet.Type = 2

# Rename Eric's task from "square" to "exponent".
$ericLine = $tr.Find("Eric: square")
$ericLine.Text = "Eric: exponent"
